$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualiza base de datos EC" - update the "Periodo Mora" (arrears period)
# column for the existing workers so the periods now read in ascending
# order (1607, 1608, 1609, 1610) instead of descending (1610, 1609, 1608,
# 1607).
$ws.Range("E16").Value = "1607"
$ws.Range("E17").Value = "1608"
$ws.Range("E18").Value = "1609"
$ws.Range("E19").Value = "1610"
